$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# This script regenerates the handback-status report for two new/renamed
# source files:
#   54668c75-76de-4c0a-95ca-b600d9003c0e.md  ->  1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md
#   9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md  ->  ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md
# along with refreshed xlf correspondence files and timestamps.
# -----------------------------------------------------------------------

# ====================== Sheet: Overview ======================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsOverview.Range("B2").Value = "e2e\1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsOverview.Range("G2").Value = "2016-08-15 20:58:41"

$wsOverview.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsOverview.Range("B3").Value = "e2e\ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsOverview.Range("G3").Value = "2016-08-15 20:58:41"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "e2e\1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "e2e\ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")

# ====================== Sheet: zh-cn ======================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsZhCn.Range("G2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 20:58:36"
$wsZhCn.Range("I2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsZhCn.Range("J2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-15 20:58:52"

$wsZhCn.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsZhCn.Range("G3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-15 20:58:36"
$wsZhCn.Range("I3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsZhCn.Range("J3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-15 20:58:52"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d5b67d0858af23a458d2467c7de1529f8f8715ab/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d5b67d0858af23a458d2467c7de1529f8f8715ab/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")

# ====================== Sheet: de-de ======================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsDeDe.Range("G2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 20:58:41"
$wsDeDe.Range("I2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md"
$wsDeDe.Range("J2").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-15 20:59:00"

$wsDeDe.Range("A3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsDeDe.Range("G3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-15 20:58:41"
$wsDeDe.Range("I3").Value = "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md"
$wsDeDe.Range("J3").Value = "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.33bb507d57f63cf80d06c5c9b063785b968b8674.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-15 20:59:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/eb3c77e157ecbade6f57ff1a9591a7e57f9a0bb0/e2e/54668c75-76de-4c0a-95ca-b600d9003c0e.md", "", "", "1f9f3455-e85d-457f-b7b3-ef852e7d7f29.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55e7f1322053d3a6fc1deef8cb824ee2679602c4/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/eb3c77e157ecbade6f57ff1a9591a7e57f9a0bb0/e2e/9f2d4457-f4c6-4bd8-a03d-1cb91a9e56ea.md", "", "", "ffff1bfafb35-02c8-49ea-8a17-273f0f6519a8.md")
